$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.256669640541077
$ws.Range("B1").Value = 2.288324594497681
$ws.Range("C1").Value = 3.907950401306152
$ws.Range("D1").Value = 2.748371839523315
$ws.Range("E1").Value = 1.364812135696411
